$d = $word.ActiveDocument

# The Pearson logo (footers) is currently saved as "image1.png" and needs to
# become "image2.png"; the BTEC logo (headers) is currently saved as
# "image2.jpg" and needs to become "image1.jpg". Walk every section's
# headers/footers and rename each inline picture accordingly.

function Rename-InlineLogo($range) {
    if ($range.InlineShapes.Count -eq 0) { return }
    for ($k = 1; $k -le $range.InlineShapes.Count; $k++) {
        $shape = $range.InlineShapes($k)
        # Selecting the shape first and renaming through the Selection's
        # InlineShapes collection is what actually persists the change.
        $shape.Select()
        $sel = $word.Selection
        $selShape = $sel.InlineShapes(1)
        if ($selShape.AlternativeText -eq "BTec_Logo-Orange") {
            $selShape.Name = "image1.jpg"
        } elseif ($selShape.AlternativeText -like "*PearsonLogo.png") {
            $selShape.Name = "image2.png"
        }
    }
}

for ($i = 1; $i -le $d.Sections.Count; $i++) {
    $sec = $d.Sections($i)

    for ($h = 1; $h -le $sec.Headers.Count; $h++) {
        $hdr = $sec.Headers($h)
        if ($hdr.Exists) {
            Rename-InlineLogo $hdr.Range
        }
    }

    for ($f = 1; $f -le $sec.Footers.Count; $f++) {
        $ftr = $sec.Footers($f)
        if ($ftr.Exists) {
            Rename-InlineLogo $ftr.Range
        }
    }
}
